$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(
    "04/10/2019",
    "05/10/2019",
    "06/10/2019",
    "07/10/2019",
    "08/10/2019",
    "09/10/2019",
    "10/10/2019",
    "11/10/2019",
    "12/10/2019",
    "13/10/2019",
    "14/10/2019",
    "15/10/2019",
    "16/10/2019",
    "17/10/2019",
    "18/10/2019",
    "19/10/2019",
    "20/10/2019",
    "21/10/2019",
    "22/10/2019",
    "23/10/2019",
    "24/10/2019",
    "25/10/2019",
    "26/10/2019",
    "27/10/2019",
    "28/10/2019",
    "29/10/2019"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$i]
}
